# Update DAMSLTag (col I) and DialogAct (col J) values for specific rows
# following a re-run of SGNN dialog act annotation after transcript cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 7; I = "sv"; J = "Statement-opinion" }
    @{ Row = 16; I = "%"; J = "Uninterpretable" }
    @{ Row = 33; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 34; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 42; I = "sv"; J = "Statement-opinion" }
    @{ Row = 47; I = "ba"; J = "Appreciation" }
    @{ Row = 52; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 71; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 73; I = "sv"; J = "Statement-opinion" }
    @{ Row = 80; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 83; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 89; I = "aa"; J = "Agree/Accept" }
    @{ Row = 91; I = "aa"; J = "Agree/Accept" }
    @{ Row = 101; I = "aa"; J = "Agree/Accept" }
    @{ Row = 102; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 104; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 108; I = "sv"; J = "Statement-opinion" }
    @{ Row = 113; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 114; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 133; I = "sv"; J = "Statement-opinion" }
    @{ Row = 135; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 142; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 145; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 147; I = "%"; J = "Uninterpretable" }
    @{ Row = 151; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 160; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 165; I = "ba"; J = "Appreciation" }
    @{ Row = 172; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 173; I = "sv"; J = "Statement-opinion" }
    @{ Row = 178; I = "ba"; J = "Appreciation" }
    @{ Row = 189; I = "sv"; J = "Statement-opinion" }
    @{ Row = 190; I = "%"; J = "Uninterpretable" }
    @{ Row = 193; I = "sv"; J = "Statement-opinion" }
    @{ Row = 194; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 196; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 199; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 223; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 230; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 241; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 245; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 246; I = "aa"; J = "Agree/Accept" }
    @{ Row = 267; I = "sv"; J = "Statement-opinion" }
    @{ Row = 291; I = "%"; J = "Uninterpretable" }
    @{ Row = 297; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 309; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 327; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 332; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 334; I = "sv"; J = "Statement-opinion" }
    @{ Row = 354; I = "sv"; J = "Statement-opinion" }
    @{ Row = 356; I = "sv"; J = "Statement-opinion" }
    @{ Row = 358; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 362; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 366; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 378; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 381; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 407; I = "sv"; J = "Statement-opinion" }
    @{ Row = 412; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 425; I = "sv"; J = "Statement-opinion" }
    @{ Row = 430; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 440; I = "sv"; J = "Statement-opinion" }
    @{ Row = 442; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 444; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 449; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 452; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 458; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.I
    $ws.Cells.Item($change.Row, 10).Value = $change.J
}
